$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "309.98"
Set-TextValue "E2" "1.29%"
Set-TextValue "D3" "35.65"
Set-TextValue "E3" "-1.43%"
Set-TextValue "D4" "5.121"
Set-TextValue "E4" "1.53%"
Set-TextValue "D5" "0.08211"
Set-TextValue "E5" "4.34%"
Set-TextValue "D6" "2.057"
Set-TextValue "E6" "-9.52%"
Set-TextValue "D7" "7.973"
Set-TextValue "E7" "-0.34%"
Set-TextValue "D8" "4.132"
Set-TextValue "E8" "-0.46%"
Set-TextValue "D9" "2.896"
Set-TextValue "E9" "8.86%"
Set-TextValue "D10" "0.9268"
Set-TextValue "E10" "-0.09%"
Set-TextValue "D11" "0.1077"
Set-TextValue "E11" "9.80%"
Set-TextValue "D12" "0.1920"
Set-TextValue "E12" "2.89%"
Set-TextValue "E13" "7.88%"
Set-TextValue "D14" "0.03601"
Set-TextValue "E14" "-4.10%"
Set-TextValue "D15" "0.09911"
Set-TextValue "E15" "-0.20%"
Set-TextValue "D16" "0.001441"
Set-TextValue "E16" "0.43%"
Set-TextValue "D17" "0.005880"
Set-TextValue "E17" "4.58%"
Set-TextValue "D18" "3.474"
Set-TextValue "E18" "0.40%"
Set-TextValue "D19" "0.3425"
Set-TextValue "E19" "1.74%"
Set-TextValue "E20" "-0.68%"
Set-TextValue "D21" "5.102"
Set-TextValue "E21" "-0.29%"
Set-TextValue "D22" "0.2193"
Set-TextValue "E22" "-2.64%"
Set-TextValue "D23" "0.04552"
Set-TextValue "E23" "-0.63%"
Set-TextValue "D24" "0.001224"
Set-TextValue "E24" "-0.91%"
Set-TextValue "D25" "0.004803"
Set-TextValue "E25" "0.59%"
Set-TextValue "E26" "-3.90%"
Set-TextValue "D27" "0.0004453"
Set-TextValue "E27" "-6.02%"
Set-TextValue "D39" "0.01978"
Set-TextValue "E39" "2.96%"
Set-TextValue "D40" "0.04914"
Set-TextValue "E40" "-0.83%"
Set-TextValue "D41" "0.007617"
Set-TextValue "E41" "-2.43%"
Set-TextValue "D42" "0.009853"
Set-TextValue "E42" "26.08%"
Set-TextValue "D43" "0.1383"
Set-TextValue "E43" "-0.53%"
Set-TextValue "D44" "0.002116"
Set-TextValue "E44" "-0.88%"
Set-TextValue "D45" "0.01155"
Set-TextValue "E45" "1.20%"
Set-TextValue "D46" "0.00006504"
Set-TextValue "E46" "5.42%"
Set-TextValue "E47" "-0.01%"
Set-TextValue "D48" "175.39"
Set-TextValue "E48" "238.80%"
Set-TextValue "E49" "-16.87%"
Set-TextValue "D50" "0.00002102"
Set-TextValue "E50" "-0.01%"
Set-TextValue "D51" "0.0002002"
Set-TextValue "E51" "-0.01%"
